# "por fin byobu keybindings" - replace the byobu sheet's tmux-era shortcut
# table with the real byobu keybindings, add a quick-reference command list
# in column D (rows 2-17), make "byobu" the active/selected sheet, and flip
# the sheet to landscape printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("byobu")

# ---- Column A/B formatting: rows 2-10 already carry the banded s13/s12
#      look; extend the same banding down through row 13 ----
$ws.Range("A8:B9").Copy() | Out-Null
$ws.Range("A10:B13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A2").Value2  = "Move to Split r/l"
$ws.Range("B2").Value2  = "Shift + right/left"
$ws.Range("A3").Value2  = "Show numbers"
$ws.Range("B3").Value2  = "Shift + F6"
$ws.Range("A4").Value2  = "Detach & log out"
$ws.Range("B4").Value2  = "F6"
$ws.Range("A5").Value2  = "Kill Window/Split"
$ws.Range("B5").Value2  = "Ctrl + D (exit RET)"
$ws.Range("A6").Value2  = "Rename Window"
$ws.Range("B6").Value2  = "F8"
$ws.Range("A7").Value2  = "Show Help"
$ws.Range("B7").Value2  = "Shift + F1"
$ws.Range("A8").Value2  = "New Window"
$ws.Range("B8").Value2  = "F2"
$ws.Range("A9").Value2  = "Move to Window"
$ws.Range("B9").Value2  = "Alt + left/right (F3-F4)"
$ws.Range("A10").Value2 = "Ctrl + F9 "
$ws.Range("B10").Value2 = "command all windows"
$ws.Range("A11").Value2 = "Shift + F9"
$ws.Range("B11").Value2 = "command all splits"
$ws.Range("A12").Value2 = "Alt + F9"
$ws.Range("B12").Value2 = "toggle type all splits"
$ws.Range("A13").Value2 = "Shift + F12 "
$ws.Range("B13").Value2 = "Toggle keybindings"

# ---- Column C: only the header row keeps a value now, un-banded (s12) ----
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$ws.Range("C2").Value2 = "F12"
$ws.Range("C3:C11").Clear() | Out-Null

# ---- Column D: new byobu quick-reference list, rows 2-17, all un-banded
#      (style s12 throughout, unlike the banded A/B columns) ----
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$ws.Range("D4:D17").PasteSpecial(-4122) | Out-Null     # xlPasteFormats

$ws.Range("D2").Value2  = "% split vert"
$ws.Range("D3").Value2  = "| Split horiz"
$ws.Range("D4").Value2  = "c  create window"
$ws.Range("D5").Value2  = "w  list windows"
$ws.Range("D6").Value2  = "n  next window"
$ws.Range("D7").Value2  = "p  previous window"
$ws.Range("D8").Value2  = "f  find window"
$ws.Range("D9").Value2  = ",  name window"
$ws.Range("D10").Value2 = "&  kill window"
$ws.Range("D11").Value2 = "o  swap panes"
$ws.Range("D12").Value2 = "q  show pane numbers"
$ws.Range("D13").Value2 = "k  kill pane"
$ws.Range("D14").Value2 = "d  detach"
$ws.Range("D15").Value2 = "t  big clock"
$ws.Range("D16").Value2 = "?  list shortcuts"
$ws.Range("D17").Value2 = ":  prompt"

# ---- Row heights: keep the sheet's standard 17.25 row height explicit on
#      the newly-used rows 11-17 (rows 2-10 already carry it) ----
for ($r = 11; $r -le 17; $r++) {
    $ws.Rows.Item($r).RowHeight = 17.25
}

# ---- Column widths (best-effort match of the autofit result) ----
$ws.Columns.Item(1).ColumnWidth = 21.45
$ws.Columns.Item(2).ColumnWidth = 29.02
$ws.Columns.Item(3).ColumnWidth = 6.59
$ws.Columns.Item(4).ColumnWidth = 30.88

# ---- Page orientation -> landscape ----
$ws.PageSetup.Orientation = 2   # xlLandscape

# ---- Selection + make byobu the active/selected sheet & tab ----
$ws.Range("B16").Select() | Out-Null
$ws.Activate()
